$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) values are numeric-looking text strings
# (e.g. "1.000", "0.06580", "  -0.83%  "). Force each target cell to Text
# format before writing so Excel keeps the exact original formatting
# instead of silently converting it to a floating point number.
foreach ($ref in @("D2","E2","D3","E3","D4","E4","D5","E5","E6","D7","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","D19","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","D29","E29","D30","E30","D31","E31","D32","E32","D33","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","D45","E45","D46","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "28.016.99"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "1.801.96"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "310.06"
$ws.Range("E5").Value = "  -1.53%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "0.5076"
$ws.Range("E7").Value = "  -3.53%  "
$ws.Range("D8").Value = "0.3857"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").Value = "0.08917"
$ws.Range("E9").Value = "  +11.03%  "
$ws.Range("D10").Value = "1.096"
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("D11").Value = "40.77"
$ws.Range("E11").Value = "  -1.78%  "
$ws.Range("D12").Value = "6.376"
$ws.Range("E12").Value = "  +0.13%  "
$ws.Range("D13").Value = "1.000"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "20.32"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "7.292"
$ws.Range("E15").Value = "  -1.14%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "1.791.22"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("B18").Value = "Litecoin"
$ws.Range("C18").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D18").Value = "92.11"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Range("D19").Value = "0.06580"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "17.26"
$ws.Range("E21").Value = "  -1.01%  "
$ws.Range("D22").Value = "6.013"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("D23").Value = "28.056.29"
$ws.Range("E23").Value = "  -0.90%  "
$ws.Range("D24").Value = "11.05"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").Value = "2.222"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").Value = "158.46"
$ws.Range("E26").Value = "  -0.94%  "
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.009.43"
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.407"
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("D29").Value = "20.30"
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("D30").Value = "127.53"
$ws.Range("E30").Value = "  +3.42%  "
$ws.Range("D31").Value = "0.1090"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "1.046"
$ws.Range("E32").Value = "  -1.63%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.562"
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "3.642"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").Value = "0.06941"
$ws.Range("E35").Value = "  -4.57%  "
$ws.Range("D36").Value = "9.033"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").Value = "0.02335"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.2169"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "5.011"
$ws.Range("E39").Value = "  -2.24%  "
$ws.Range("D40").Value = "11.41"
$ws.Range("E40").Value = "  -8.39%  "
$ws.Range("D41").Value = "0.6123"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").Value = "1.0000"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "1.151"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("D44").Value = "13.28"
$ws.Range("E44").Value = "  +0.31%  "
$ws.Range("D45").Value = "1.294"
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("D46").Value = "0.5898"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "3.707"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "125.25"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "1.929"
$ws.Range("E49").Value = "  -0.55%  "
$ws.Range("D50").Value = "1.184"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("D51").Value = "0.06740"
$ws.Range("E51").Value = "  -1.58%  "

Write-Host "Applied 115 cell updates"
